{"js": "// Update the 25 three-digit-by-one-digit multiplication problems/answers\n// in the practice worksheet table to the new set of problems.\nconst replacements = [\n  [\"953\u00d77=6671\", \"494\u00d79=4446\"],\n  [\"409\u00d73=1227\", \"335\u00d74=1340\"],\n  [\"144\u00d72=288\", \"734\u00d77=5138\"],\n  [\"849\u00d74=3396\", \"179\u00d75=895\"],\n  [\"398\u00d73=1194\", \"190\u00d76=1140\"],\n  [\"181\u00d72=362\", \"428\u00d78=3424\"],\n  [\"220\u00d73=660\", \"252\u00d75=1260\"],\n  [\"203\u00d74=812\", \"514\u00d79=4626\"],\n  [\"401\u00d79=3609\", \"901\u00d77=6307\"],\n  [\"155\u00d74=620\", \"311\u00d72=622\"],\n  [\"308\u00d77=2156\", \"772\u00d76=4632\"],\n  [\"482\u00d75=2410\", \"686\u00d78=5488\"],\n  [\"477\u00d76=2862\", \"622\u00d76=3732\"],\n  [\"322\u00d73=966\", \"444\u00d78=3552\"],\n  [\"224\u00d74=896\", \"760\u00d76=4560\"],\n  [\"897\u00d77=6279\", \"930\u00d77=6510\"],\n  [\"306\u00d72=612\", \"317\u00d75=1585\"],\n  [\"810\u00d77=5670\", \"982\u00d75=4910\"],\n  [\"282\u00d75=1410\", \"863\u00d78=6904\"],\n  [\"508\u00d72=1016\", \"806\u00d77=5642\"],\n  [\"836\u00d77=5852\", \"431\u00d79=3879\"],\n  [\"352\u00d79=3168\", \"546\u00d78=4368\"],\n  [\"383\u00d74=1532\", \"480\u00d74=1920\"],\n  [\"236\u00d75=1180\", \"714\u00d76=4284\"],\n  [\"369\u00d75=1845\", \"201\u00d77=1407\"]\n];\n\nconst body = context.document.body;\n\n// Each \"old\" string is unique in the document, so a plain text search\n// (matchCase true, no wildcards) safely targets exactly one run each time.\n// Replacing in-place on the found Range preserves the original run\n// formatting (font / size) instead of re-creating a new run.\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 three-digit-by-one-digit multiplication problems/answers\n# in the practice worksheet table to the new set of problems.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"953\u00d77=6671\", \"494\u00d79=4446\"),\n  @(\"409\u00d73=1227\", \"335\u00d74=1340\"),\n  @(\"144\u00d72=288\", \"734\u00d77=5138\"),\n  @(\"849\u00d74=3396\", \"179\u00d75=895\"),\n  @(\"398\u00d73=1194\", \"190\u00d76=1140\"),\n  @(\"181\u00d72=362\", \"428\u00d78=3424\"),\n  @(\"220\u00d73=660\", \"252\u00d75=1260\"),\n  @(\"203\u00d74=812\", \"514\u00d79=4626\"),\n  @(\"401\u00d79=3609\", \"901\u00d77=6307\"),\n  @(\"155\u00d74=620\", \"311\u00d72=622\"),\n  @(\"308\u00d77=2156\", \"772\u00d76=4632\"),\n  @(\"482\u00d75=2410\", \"686\u00d78=5488\"),\n  @(\"477\u00d76=2862\", \"622\u00d76=3732\"),\n  @(\"322\u00d73=966\", \"444\u00d78=3552\"),\n  @(\"224\u00d74=896\", \"760\u00d76=4560\"),\n  @(\"897\u00d77=6279\", \"930\u00d77=6510\"),\n  @(\"306\u00d72=612\", \"317\u00d75=1585\"),\n  @(\"810\u00d77=5670\", \"982\u00d75=4910\"),\n  @(\"282\u00d75=1410\", \"863\u00d78=6904\"),\n  @(\"508\u00d72=1016\", \"806\u00d77=5642\"),\n  @(\"836\u00d77=5852\", \"431\u00d79=3879\"),\n  @(\"352\u00d79=3168\", \"546\u00d78=4368\"),\n  @(\"383\u00d74=1532\", \"480\u00d74=1920\"),\n  @(\"236\u00d75=1180\", \"714\u00d76=4284\"),\n  @(\"369\u00d75=1845\", \"201\u00d77=1407\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n\n  # wdFindContinue = 1, wdReplaceAll = 2\n  # Each \"old\" string is unique within the document, so MatchCase/MatchWholeWord\n  # ensure exactly one run is targeted and its formatting is preserved in place.\n  $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
